$d = $word.ActiveDocument

# Rename Model.Tasks -> Model.Objectives throughout the report template
# (ReportModel: renamed Task to Objectives for consistency)

$d.Content.Find.Execute(
    "foreach (var task in Model.Tasks) { ", $true, $false, $false, $false, $false,
    $true, 1, $false, "foreach (var task in Model.Objectives) { ", 2) | Out-Null

$d.Content.Find.Execute(
    "Model.Tasks.Last()", $true, $false, $false, $false, $false,
    $true, 1, $false, "Model.Objectives.Last()", 2) | Out-Null
